# Auto-generated Excel COM-interop edit script
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row height adjustments ---
$ws.Rows.Item(2).RowHeight = 27.75
$ws.Rows.Item(3).RowHeight = 27.75
$ws.Rows.Item(4).RowHeight = 27.75
$ws.Rows.Item(5).RowHeight = 27.75
$ws.Rows.Item(6).RowHeight = 27.75
$ws.Rows.Item(7).RowHeight = 27.75
$ws.Rows.Item(8).RowHeight = 27.75
$ws.Rows.Item(9).RowHeight = 27.75
$ws.Rows.Item(10).RowHeight = 27.75
$ws.Rows.Item(11).RowHeight = 41.25
$ws.Rows.Item(12).RowHeight = 41.25
$ws.Rows.Item(13).RowHeight = 27.75
$ws.Rows.Item(14).RowHeight = 27.75
$ws.Rows.Item(15).RowHeight = 27.75
$ws.Rows.Item(16).RowHeight = 27.75
$ws.Rows.Item(17).RowHeight = 27.75
$ws.Rows.Item(18).RowHeight = 68.25
$ws.Rows.Item(25).RowHeight = 27.75
$ws.Rows.Item(43).RowHeight = 27.75
$ws.Rows.Item(44).RowHeight = 27.75
$ws.Rows.Item(45).RowHeight = 27.75
$ws.Rows.Item(52).RowHeight = 27.75

# --- Vertical alignment: bottom -> top for F/H/I remarks & pending-work columns (rows 19-60) ---
$rngTop = $ws.Range('F19,H19,I19,F20,H20,I20,F21,H21,I21,F22,H22,I22,F23,H23,I23,F24,H24,I24,F25,H25,I25,F26,H26,I26,F27,H27,I27,F28,H28,I28,F29,H29,I29,F30,H30,I30,F31,H31,I31,F32,H32,I32,F33,H33,I33,F34,H34,I34,F35,H35,I35,F36,H36,I36,F37,H37,I37,F38,H38,I38,F39,H39,I39,F40,H40,I40,F41,H41,I41,F42,H42,I42,F43,H43,I43,F44,H44,I44,F45,H45,I45,F46,H46,I46,F47,H47,I47,F48,H48,I48,F49,H49,I49,F50,H50,I50,F51,H51,I51,F52,H52,I52,F53,H53,I53,F54,H54,I54,F55,H55,I55,F56,H56,I56,F57,H57,I57,F58,H58,I58,F59,H59,I59,F60,H60,I60')
$rngTop.VerticalAlignment = -4160  # xlVAlignTop

# --- Update existing rows 61-66: fill in F/H/I content and flip G from pending to implemented ---
$ws.Range("F61").Value = 'Analytics page now loads summary metrics via /api/analytics/summary and supports basic filtering by strategy and date range.'
$ws.Range("G61").Value = 'implemented'
$ws.Range("H61").Value = 'Summary card shows trades, total P&L, win rate, avg win/loss, and max drawdown; a Rebuild button triggers backend rebuild-trades and refreshes metrics.'
$ws.Range("I61").Value = 'Expose more advanced metrics (e.g., expectancy, risk-adjusted returns) once enough live data is available.'
$ws.Range("F62").Value = 'Added lightweight SVG-based cumulative P&L line chart and P&L by symbol bar chart driven by analytics trades.'
$ws.Range("G62").Value = 'implemented'
$ws.Range("H62").Value = 'Charts update automatically when filters change and are kept deliberately simple without adding new dependencies.'
$ws.Range("I62").Value = 'Consider richer charting (e.g., per-day P&L) if/when a chart library is introduced.'
$ws.Range("F63").Value = 'Trades table shows closed trades with IST timestamps, strategy, symbol, and colored P&L, using filters for strategy and date range from the Analytics header.'
$ws.Range("G63").Value = 'implemented'
$ws.Range("H63").Value = 'Gives a clear audit-style view of trades aligned with the summary and charts.'
$ws.Range("I63").Value = 'Add more filters (e.g., symbol, min/max P&L) if needed for heavier analysis.'
$ws.Range("F64").Value = 'Added JSON-formatted structured logging with correlation IDs via a FastAPI middleware; key flows (webhook ingestion, order execution, Zerodha connect) now emit structured events.'
$ws.Range("G64").Value = 'implemented'
$ws.Range("H64").Value = 'Logs include correlation_id so a single alert/order can be traced across webhook, risk, and broker calls.'
$ws.Range("I64").Value = 'Consider integrating with an external log aggregator if the app is deployed beyond local single-user usage.'
$ws.Range("F64,H64,I64").VerticalAlignment = -4160
$ws.Range("F64,H64,I64").WrapText = $false
$ws.Range("F65").Value = 'Standardized several error paths (risk rejection, Zerodha order failures) and ensured correlation IDs are included in both logs and HTTP responses.'
$ws.Range("G65").Value = 'implemented'
$ws.Range("H65").Value = 'Critical errors now have consistent HTTP status codes and messages while being traceable via X-Request-ID.'
$ws.Range("I65").Value = 'Extend normalization to any remaining ad-hoc error paths as they are discovered.'
$ws.Range("F65,H65,I65").VerticalAlignment = -4160
$ws.Range("F65,H65,I65").WrapText = $false
$ws.Range("F66").Value = 'Added a simple System Events page that surfaces recent client-side errors/warnings captured in memory during the current session.'
$ws.Range("G66").Value = 'implemented'
$ws.Range("H66").Value = 'Gives a quick at-a-glance view of recent issues without setting up a full log viewer.'
$ws.Range("I66").Value = 'If needed, extend this view to pull server-side logs or important alerts once a log backend is introduced.'
$ws.Range("F66,H66,I66").VerticalAlignment = -4160
$ws.Range("F66,H66,I66").WrapText = $false

# --- Add new rows 72-74 (S08/G04 tasks) ---
$ws.Rows.Item(72).RowHeight = 41.75
$ws.Range("A72").Value = 'S08'
$ws.Range("B72").Value = 'G04'
$ws.Range("C72").Value = 'Backend log aggregation and server log UI'
$ws.Range("D72").Value = 'S08_G04_TB001'
$ws.Range("E72").Value = 'Persist important backend events (alerts, orders, broker, risk) in a system_events table.'
$ws.Range("F72").Value = 'Introduced a system_events table and helpers so key backend events (alerts ingested, orders executed, Zerodha connects, risk rejections) are persisted.'
$ws.Range("G72").Value = 'implemented'
$ws.Range("H72").Value = 'Important events can now be queried and shown in the UI without parsing textual logs.'
$ws.Range("I72").Value = 'Extend coverage to any additional event categories we decide to track in future (e.g., background jobs, sync failures).'
$ws.Range("A72:E72").VerticalAlignment = -4160
$ws.Range("A72:E72").WrapText = $false
$ws.Range("G72").VerticalAlignment = -4160
$ws.Range("G72").WrapText = $false
$ws.Range("F72,H72,I72").VerticalAlignment = -4160
$ws.Range("F72,H72,I72").WrapText = $true
$ws.Rows.Item(73).RowHeight = 41.75
$ws.Range("A73").Value = 'S08'
$ws.Range("B73").Value = 'G04'
$ws.Range("C73").Value = 'Backend log aggregation and server log UI'
$ws.Range("D73").Value = 'S08_G04_TB002'
$ws.Range("E73").Value = 'Expose API endpoints to query recent system events for the UI.'
$ws.Range("F73").Value = 'Added /api/system-events/ endpoint with filters for level, category, and limit, returning recent events in reverse chronological order.'
$ws.Range("G73").Value = 'implemented'
$ws.Range("H73").Value = 'Provides a simple API surface that the UI and tools can use to inspect important backend events.'
$ws.Range("I73").Value = 'Consider pagination or time-based filters if the event volume grows significantly.'
$ws.Range("A73:E73").VerticalAlignment = -4160
$ws.Range("A73:E73").WrapText = $false
$ws.Range("G73").VerticalAlignment = -4160
$ws.Range("G73").WrapText = $false
$ws.Range("F73,H73,I73").VerticalAlignment = -4160
$ws.Range("F73,H73,I73").WrapText = $true
$ws.Rows.Item(74).RowHeight = 41.75
$ws.Range("A74").Value = 'S08'
$ws.Range("B74").Value = 'G04'
$ws.Range("C74").Value = 'Backend log aggregation and server log UI'
$ws.Range("D74").Value = 'S08_G04_TF003'
$ws.Range("E74").Value = 'Show a table of recent backend events (alerts, orders, broker events) in the web app.'
$ws.Range("F74").Value = 'System Events page now displays recent backend events from the system_events table along with client-side events for the current browser session.'
$ws.Range("G74").Value = 'implemented'
$ws.Range("H74").Value = 'Gives an at-a-glance view of both server and client issues without leaving the app.'
$ws.Range("I74").Value = 'Enhance filtering (by level/category/time) in the UI if needed during real operations.'
$ws.Range("A74:E74").VerticalAlignment = -4160
$ws.Range("A74:E74").WrapText = $false
$ws.Range("G74").VerticalAlignment = -4160
$ws.Range("G74").WrapText = $false
$ws.Range("F74,H74,I74").VerticalAlignment = -4160
$ws.Range("F74,H74,I74").WrapText = $true

# --- Update sheet view: scroll position and selection (best effort) ---
$win = $excel.ActiveWindow
$win.ScrollRow = 53
$win.ScrollColumn = 1
$ws.Range("1:1048576").Select()
